$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")  # the "Metadata" property/value sheet

# Bump the published version and republish date
$ws.Range("B3").Value = "0.1.1"
$ws.Range("B8").Value = "2023-06-02T12:02:38+02:00"

# The extension now declares an additional usage context. The existing
# "Context" row keeps its label but its value becomes the new context,
# and the previous context value moves down into a new row right below it.
$previousContext = $ws.Range("B20").Text
$ws.Range("B20").Value = "element:ContactPoint"

# Duplicate row 20's formatting onto the new row 21 before filling it in.
$ws.Range("A20:B20").Copy()
$ws.Range("A21:B21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A21").Value = "Context"
$ws.Range("B21").Value = $previousContext
